# Fruta / hortaliza, semanal
# Insert a new weekly record at row 203 (pushing the existing history
# down by one row) on the single data sheet of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 203:254 down to 204:255, creating a blank row 203.
$ws.Rows.Item(203).Insert()

# Populate the new row with the latest observation.
$ws.Cells.Item(203, 1).Value  = 10
$ws.Cells.Item(203, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(203, 3).Value  = "La Araucanía"
$ws.Cells.Item(203, 4).Value  = 44736
$ws.Cells.Item(203, 5).Value  = 9
$ws.Cells.Item(203, 6).Value  = 100112039
$ws.Cells.Item(203, 7).Value  = "Ciboulette"
$ws.Cells.Item(203, 8).Value  = "Sin especificar"
$ws.Cells.Item(203, 9).Value  = "Primera"
$ws.Cells.Item(203, 10).Value = 30
$ws.Cells.Item(203, 11).Value = 6000
$ws.Cells.Item(203, 12).Value = 6000
$ws.Cells.Item(203, 13).Value = 6000
$ws.Cells.Item(203, 14).Value = "$/docena de atados"
$ws.Cells.Item(203, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(203, 16).Value = 2000
$ws.Cells.Item(203, 17).Value = 3
$ws.Cells.Item(203, 18).Value = "Hortaliza"
